$wb = $excel.ActiveWorkbook
$agg = $wb.Worksheets.Item("Agg")

$headerFill = 16764057    # RGB(153,204,255) = 0x99CCFF -> BGR-ish value for COM Color (R + G*256 + B*65536)
$moneyFormat = '"$"#,##0.00'
$colWidthOffset = 0.8333333333

function New-AggSheet {
    param($prevSheet, $sheetName)
    $ws = $wb.Worksheets.Add($null, $prevSheet)
    $ws.Name = $sheetName
    return $ws
}

function Format-Header {
    param($rng)
    $rng.Borders.LineStyle = 1
    $rng.Font.Bold = $true
    $rng.Interior.Color = $headerFill
}

function Format-Plain {
    param($rng)
    $rng.Borders.LineStyle = 1
}

function Format-Money {
    param($rng)
    $rng.Borders.LineStyle = 1
    $rng.NumberFormat = $moneyFormat
}

# ---------------------------------------------------------------------------
# Sheet: Msd
# ---------------------------------------------------------------------------
$msd = New-AggSheet $agg "Msd"
$msd.Columns.Item(1).ColumnWidth = 26.42578125 - $colWidthOffset
$msd.Columns.Item(2).ColumnWidth = 12.28515625 - $colWidthOffset

$msd.Range("A1").Value = "Is A Manager"
$msd.Range("B1").Value = "Total Salary"
Format-Header $msd.Range("A1:B1")

$msd.Range("A2").Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager()" useMsd="true"><jt:forEach items="${values}" var="value">${value.getPropertyValue(0)}'
Format-Plain $msd.Range("A2")

$msd.Range("B2").Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'
Format-Money $msd.Range("B2")

$msd.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Sheet: Rollup
# ---------------------------------------------------------------------------
$rollup = New-AggSheet $msd "Rollup"
$rollup.Columns.Item(1).ColumnWidth = 26.42578125 - $colWidthOffset
$rollup.Columns.Item(2).ColumnWidth = 26.42578125 - $colWidthOffset
$rollup.Columns.Item(3).ColumnWidth = 12.28515625 - $colWidthOffset

$rollup.Range("A1").Value = "Is A Manager"
$rollup.Range("B1").Value = "Title"
$rollup.Range("C1").Value = "Total Salary"
Format-Header $rollup.Range("A1:C1")

$rollup.Range("A2").Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title" rollup="${[0, 1]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
Format-Plain $rollup.Range("A2")

$rollup.Range("B2").Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
Format-Plain $rollup.Range("B2")

$rollup.Range("C2").Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'
Format-Money $rollup.Range("C2")

# ---------------------------------------------------------------------------
# Sheet: Rollups
# ---------------------------------------------------------------------------
$rollups = New-AggSheet $rollup "Rollups"
$rollups.Columns.Item(1).ColumnWidth = 26.42578125 - $colWidthOffset
$rollups.Columns.Item(2).ColumnWidth = 26.42578125 - $colWidthOffset
$rollups.Columns.Item(3).ColumnWidth = 31 - $colWidthOffset
$rollups.Columns.Item(4).ColumnWidth = 12.28515625 - $colWidthOffset

$rollups.Range("A1").Value = "Is A Manager"
$rollups.Range("B1").Value = "Title"
$rollups.Range("C1").Value = "Catch Phrase"
$rollups.Range("D1").Value = "Total Salary"
Format-Header $rollups.Range("A1:D1")

$rollups.Range("A2").Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title;catchPhrase" rollups="${[[1], [2]]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1);getPropertyValue(2)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
Format-Plain $rollups.Range("A2")

$rollups.Range("B2").Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
Format-Plain $rollups.Range("B2")

$rollups.Range("C2").Value = '${value.isGrouping(2) ? ''All Values'' : value.getPropertyValue(2)}'
Format-Plain $rollups.Range("C2")

$rollups.Range("D2").Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'
Format-Money $rollups.Range("D2")

# ---------------------------------------------------------------------------
# Sheet: Cube
# ---------------------------------------------------------------------------
$cube = New-AggSheet $rollups "Cube"
$cube.Columns.Item(1).ColumnWidth = 26.42578125 - $colWidthOffset
$cube.Columns.Item(2).ColumnWidth = 26.42578125 - $colWidthOffset
$cube.Columns.Item(3).ColumnWidth = 31 - $colWidthOffset
$cube.Columns.Item(4).ColumnWidth = 12.28515625 - $colWidthOffset

$cube.Range("A1").Value = "Is A Manager"
$cube.Range("B1").Value = "Title"
$cube.Range("C1").Value = "Catch Phrase"
$cube.Range("D1").Value = "Total Salary"
Format-Header $cube.Range("A1:D1")

$cube.Range("A2").Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title;catchPhrase" cube="${[0, 1, 2]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1);getPropertyValue(2)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
Format-Plain $cube.Range("A2")

$cube.Range("B2").Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
Format-Plain $cube.Range("B2")

$cube.Range("C2").Value = '${value.isGrouping(2) ? ''All Values'' : value.getPropertyValue(2)}'
Format-Plain $cube.Range("C2")

$cube.Range("D2").Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'
Format-Money $cube.Range("D2")

# ---------------------------------------------------------------------------
# Sheet: GroupingSets
# ---------------------------------------------------------------------------
$gsets = New-AggSheet $cube "GroupingSets"
$gsets.Columns.Item(1).ColumnWidth = 26.42578125 - $colWidthOffset
$gsets.Columns.Item(2).ColumnWidth = 26.42578125 - $colWidthOffset
$gsets.Columns.Item(3).ColumnWidth = 31 - $colWidthOffset
$gsets.Columns.Item(4).ColumnWidth = 12.28515625 - $colWidthOffset

$gsets.Range("A1").Value = "Is A Manager"
$gsets.Range("B1").Value = "Title"
$gsets.Range("C1").Value = "Catch Phrase"
$gsets.Range("D1").Value = "Total Salary"
Format-Header $gsets.Range("A1:D1")

$gsets.Range("A2").Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title;catchPhrase" groupingSets="${[[0], [1, 2]]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1);getPropertyValue(2)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
Format-Plain $gsets.Range("A2")

$gsets.Range("B2").Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
Format-Plain $gsets.Range("B2")

$gsets.Range("C2").Value = '${value.isGrouping(2) ? ''All Values'' : value.getPropertyValue(2)}'
Format-Plain $gsets.Range("C2")

$gsets.Range("D2").Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'
Format-Money $gsets.Range("D2")

# Restore the originally active sheet.
$agg.Activate()
